$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 129 (pushes existing rows 129-206 down to 131-208)
$ws.Range("A129:A130").EntireRow.Insert()

# Fill in the new row 129 with the new weekly record
$ws.Range("A129").Value = 10
$ws.Range("B129").Value = "Vega Modelo de Temuco"
$ws.Range("C129").Value = "La Araucanía"
$ws.Range("D129").Value = 44767
$ws.Range("E129").Value = 9
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100104
$ws.Range("H129").Value = "Frutos de pepita"
$ws.Range("I129").Value = 100104003
$ws.Range("J129").Value = "Membrillo"
$ws.Range("K129").Value = "Champion"
$ws.Range("L129").Value = "Especial"
$ws.Range("M129").Value = 8
$ws.Range("N129").Value = 250000
$ws.Range("O129").Value = 250000
$ws.Range("P129").Value = 250000
$ws.Range("Q129").Value = "`$/bins (450 kilos)"
$ws.Range("R129").Value = "Región de O'Higgins"
$ws.Range("S129").Value = 556
$ws.Range("T129").Value = 450

# Fill in the new row 130 with the new weekly record
$ws.Range("A130").Value = 10
$ws.Range("B130").Value = "Vega Modelo de Temuco"
$ws.Range("C130").Value = "La Araucanía"
$ws.Range("D130").Value = 44767
$ws.Range("E130").Value = 9
$ws.Range("F130").Value = "Fruta"
$ws.Range("G130").Value = 100104
$ws.Range("H130").Value = "Frutos de pepita"
$ws.Range("I130").Value = 100104003
$ws.Range("J130").Value = "Membrillo"
$ws.Range("K130").Value = "Champion"
$ws.Range("L130").Value = "Primera"
$ws.Range("M130").Value = 85
$ws.Range("N130").Value = 10000
$ws.Range("O130").Value = 10000
$ws.Range("P130").Value = 10000
$ws.Range("Q130").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R130").Value = "Región de O'Higgins"
$ws.Range("S130").Value = 556
$ws.Range("T130").Value = 18
